$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format number-like text cells as Text so Excel keeps the exact
# string (with separators / trailing zeros) instead of coercing to a number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.517.32"
$ws.Range("D3").Value = "1.873.80"
$ws.Range("E4").Value = "  +0.94%  "
$ws.Range("D5").Value = "313.38"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("D7").Value = "0.4782"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("D8").Value = "0.3776"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "0.07374"
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("D10").Value = "0.9378"
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("D11").Value = "20.74"
$ws.Range("E11").Value = "  +5.76%  "
$ws.Range("D12").Value = "0.07851"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").Value = "1.906.06"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").Value = "6.589"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").Value = "90.94"
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").Value = "1.016"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "0.000008929"
$ws.Range("E18").Value = "  +3.45%  "
$ws.Range("D20").Value = "14.91"
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("D21").Value = "27.582.43"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").Value = "5.137"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "1.962"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("D25").Value = "153.81"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("D27").Value = "2.020"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").Value = "116.01"
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("D29").Value = "5.000"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").Value = "1.217"
$ws.Range("E32").Value = "  +3.74%  "
$ws.Range("D33").Value = "4.621"
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("D34").Value = "0.7499"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "2.688"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("D36").Value = "0.02070"
$ws.Range("E36").Value = "  +6.37%  "
$ws.Range("D37").Value = "1.119"
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("D38").Value = "0.05302"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "3.008"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("D40").Value = "0.5348"
$ws.Range("E40").Value = "  +2.51%  "
$ws.Range("D41").Value = "7.089"
$ws.Range("E41").Value = "  +2.45%  "
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").Value = "8.422"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").Value = "10.64"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("D45").Value = "0.4832"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("E47").Value = "  +3.55%  "
$ws.Range("D48").Value = "103.15"
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").Value = "67.41"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("D50").Value = "0.06109"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").Value = "0.8988"
$ws.Range("E51").Value = "  +1.68%  "
